$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3724414706230164
$ws.Range("B1").Value = 0.261360377073288
$ws.Range("C1").Value = 0.4038999378681183
$ws.Range("D1").Value = 4.334521293640137
$ws.Range("E1").Value = 2.3330237865448
